$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 28150.334
$ws.Range("J3").Value = 28150.334
$ws.Range("L3").Value = 28150.334
$ws.Range("N3").Value = -28378.334
$ws.Range("H58").Value = 16028
$ws.Range("I58").Value = 70
$ws.Range("J58").Value = 26666.666
$ws.Range("K58").Value = 210
$ws.Range("L58").Value = 79999.99800000001
$ws.Range("M58").Value = -60
$ws.Range("N58").Value = -80299.99800000001
$ws.Range("H62").Value = 125024136
$ws.Range("I62").Value = 333335000
$ws.Range("K62").Value = 333335000
$ws.Range("M62").Value = -333334376
$ws.Range("H64").Value = 333339330
$ws.Range("I64").Value = 1000000000
$ws.Range("K64").Value = 1000000000
$ws.Range("M64").Value = -999999752
$ws.Range("H65").Value = 125024136
$ws.Range("I65").Value = 333335000
$ws.Range("K65").Value = 1666675000
$ws.Range("M65").Value = -1666671880
$ws.Range("H67").Value = 333339330
$ws.Range("I67").Value = 1000000000
$ws.Range("K67").Value = 1000000000
$ws.Range("M67").Value = -999999142
$ws.Range("J86").Value = 7941645.5
$ws.Range("L86").Value = 7941645.5
$ws.Range("N86").Value = -7943891.5
$ws.Range("J89").Value = 7941645.5
$ws.Range("L89").Value = 39708227.5
$ws.Range("N89").Value = -39719459.5
$ws.Range("H97").Value = 14293773
$ws.Range("J97").Value = 14293773
$ws.Range("L97").Value = 42881319
$ws.Range("N97").Value = -42882311
$ws.Range("H102").Value = 28150.334
$ws.Range("J102").Value = 28150.334
$ws.Range("L102").Value = 28150.334
$ws.Range("N102").Value = -34640.334
$ws.Range("H106").Value = 76924190
$ws.Range("I106").Value = 76924190
$ws.Range("K106").Value = 76924190
$ws.Range("M106").Value = -76923559
$ws.Range("H112").Value = 4915.5
$ws.Range("J112").Value = 5079.854
$ws.Range("L112").Value = 15239.562
$ws.Range("N112").Value = -17455.562
$ws.Range("H138").Value = 3965.5115
$ws.Range("I138").Value = 1142.037
$ws.Range("J138").Value = 5215.246
$ws.Range("K138").Value = 3426.111
$ws.Range("L138").Value = 15645.738
$ws.Range("M138").Value = 1713.889
$ws.Range("N138").Value = -25925.738
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10001
$ws.Range("I45").Value = 10001
$ws.Range("K45").Value = 10001
$ws.Range("M45").Value = -9624
$ws.Range("H61").Value = 5210.8
$ws.Range("I61").Value = 2940.2917
$ws.Range("K61").Value = 2940.2917
$ws.Range("M61").Value = -2728.2917
$ws.Range("H74").Value = 1826.2069
$ws.Range("I74").Value = 842.35
$ws.Range("K74").Value = 842.35
$ws.Range("M74").Value = 31.64999999999998
$ws.Range("H77").Value = 1826.2069
$ws.Range("I77").Value = 842.35
$ws.Range("K77").Value = 4211.75
$ws.Range("M77").Value = 156.25
$ws.Range("H122").Value = 6018.615
$ws.Range("I122").Value = 2901
$ws.Range("K122").Value = 8703
$ws.Range("M122").Value = -6253
$ws.Range("H132").Value = 4851.973
$ws.Range("I132").Value = 1384.2632
$ws.Range("J132").Value = 8512.333000000001
$ws.Range("K132").Value = 4152.7896
$ws.Range("L132").Value = 25536.999
$ws.Range("M132").Value = -1622.7896
$ws.Range("N132").Value = -30596.999
$ws.Range("H136").Value = 5210.8
$ws.Range("I136").Value = 2940.2917
$ws.Range("K136").Value = 8820.875100000001
$ws.Range("M136").Value = -6270.875100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 62505544
$ws.Range("I86").Value = 25004868
$ws.Range("J86").Value = 100006216
$ws.Range("K86").Value = 25004868
$ws.Range("L86").Value = 100006216
$ws.Range("M86").Value = -25003745
$ws.Range("N86").Value = -100008462
$ws.Range("H89").Value = 62505544
$ws.Range("I89").Value = 25004868
$ws.Range("J89").Value = 100006216
$ws.Range("K89").Value = 125024340
$ws.Range("L89").Value = 500031080
$ws.Range("M89").Value = -125018724
$ws.Range("N89").Value = -500042312
$ws.Range("H94").Value = 2426.9524
$ws.Range("I94").Value = 699.3077
$ws.Range("J94").Value = 5234.375
$ws.Range("K94").Value = 699.3077
$ws.Range("L94").Value = 5234.375
$ws.Range("M94").Value = -248.3077
$ws.Range("N94").Value = -6136.375
$ws.Range("H105").Value = 2940.8823
$ws.Range("I105").Value = 2077.3333
$ws.Range("J105").Value = 3912.375
$ws.Range("K105").Value = 2077.3333
$ws.Range("L105").Value = 3912.375
$ws.Range("M105").Value = -330.3332999999998
$ws.Range("N105").Value = -7406.375
$ws.Range("H134").Value = 6131.6816
$ws.Range("I134").Value = 2192.9
$ws.Range("K134").Value = 6578.700000000001
$ws.Range("M134").Value = -4043.700000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 50.18182
$ws.Range("I7").Value = 51.285713
$ws.Range("K7").Value = 51.285713
$ws.Range("M7").Value = 61.714287
$ws.Range("H31").Value = 5535.0547
$ws.Range("I31").Value = 2347
$ws.Range("K31").Value = 2347
$ws.Range("M31").Value = -2052
$ws.Range("H34").Value = 5535.0547
$ws.Range("I34").Value = 2347
$ws.Range("K34").Value = 2347
$ws.Range("M34").Value = -2145
$ws.Range("H43").Value = 31155.666
$ws.Range("J43").Value = 31155.666
$ws.Range("L43").Value = 31155.666
$ws.Range("N43").Value = -31523.666
$ws.Range("H76").Value = 4955.4443
$ws.Range("I76").Value = 4955.4443
$ws.Range("K76").Value = 4955.4443
$ws.Range("M76").Value = -4640.4443
$ws.Range("H79").Value = 4955.4443
$ws.Range("I79").Value = 4955.4443
$ws.Range("K79").Value = 4955.4443
$ws.Range("M79").Value = -3863.4443
$ws.Range("H101").Value = 31155.666
$ws.Range("J101").Value = 31155.666
$ws.Range("L101").Value = 31155.666
$ws.Range("N101").Value = -37645.666
$ws.Range("H105").Value = 4766434.5
$ws.Range("I105").Value = 7144301
$ws.Range("K105").Value = 7144301
$ws.Range("M105").Value = -7142554
$ws.Range("H107").Value = 1945.6923
$ws.Range("J107").Value = 2495.6667
$ws.Range("L107").Value = 2495.6667
$ws.Range("N107").Value = -6335.6667
$ws.Range("H122").Value = 3033
$ws.Range("I122").Value = 2853.6428
$ws.Range("K122").Value = 8560.928400000001
$ws.Range("M122").Value = -6110.928400000001
$ws.Range("H134").Value = 3702.5
$ws.Range("I134").Value = 1532.5518
$ws.Range("K134").Value = 4597.6554
$ws.Range("M134").Value = -2062.6554
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 200000820
$ws.Range("I36").Value = 250000640
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 750001920
$ws.Range("L36").Value = 4500
$ws.Range("M36").Value = -750001751
$ws.Range("N36").Value = -4838
$ws.Range("H55").Value = 22812926
$ws.Range("I55").Value = 111111930
$ws.Range("J55").Value = 6256861.5
$ws.Range("K55").Value = 333335790
$ws.Range("L55").Value = 18770584.5
$ws.Range("M55").Value = -333335613
$ws.Range("N55").Value = -18770938.5
$ws.Range("H61").Value = 207.9
$ws.Range("I61").Value = 92.14286
$ws.Range("J61").Value = 478
$ws.Range("K61").Value = 276.42858
$ws.Range("L61").Value = 1434
$ws.Range("M61").Value = -61.42858000000001
$ws.Range("N61").Value = -1864
$ws.Range("H92").Value = 6994762.5
$ws.Range("J92").Value = 8548753
$ws.Range("L92").Value = 25646259
$ws.Range("N92").Value = -25648755
$ws.Range("H131").Value = 1546.7693
$ws.Range("I131").Value = 1345.8889
$ws.Range("J131").Value = 1998.75
$ws.Range("K131").Value = 4037.6667
$ws.Range("L131").Value = 5996.25
$ws.Range("M131").Value = 1002.3333
$ws.Range("N131").Value = -16076.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1718.3334
$ws.Range("I97").Value = 1674.75
$ws.Range("J97").Value = 1857.8
$ws.Range("K97").Value = 1674.75
$ws.Range("L97").Value = 1857.8
$ws.Range("M97").Value = -1178.75
$ws.Range("N97").Value = -2849.8
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5650.6523
$ws.Range("I7").Value = 3321.2307
$ws.Range("J7").Value = 8678.9
$ws.Range("K7").Value = 3321.2307
$ws.Range("L7").Value = 8678.9
$ws.Range("M7").Value = -3209.2307
$ws.Range("N7").Value = -8902.9
$ws.Range("H55").Value = 41667040
$ws.Range("I55").Value = 125000100
$ws.Range("K55").Value = 125000100
$ws.Range("M55").Value = -124999927
$ws.Range("H93").Value = 7911
$ws.Range("I93").Value = 5563.4287
$ws.Range("K93").Value = 5563.4287
$ws.Range("M93").Value = -4315.4287
$ws.Range("H126").Value = 5650.6523
$ws.Range("I126").Value = 3321.2307
$ws.Range("J126").Value = 8678.9
$ws.Range("K126").Value = 9963.6921
$ws.Range("L126").Value = 26036.7
$ws.Range("M126").Value = -7493.6921
$ws.Range("N126").Value = -30976.7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 200000000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = ""
$ws.Range("H84").Value = 200000000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = ""
$ws.Range("H122").Value = 2899.641
$ws.Range("I122").Value = 1439.625
$ws.Range("K122").Value = 4318.875
$ws.Range("M122").Value = -1868.875
$ws.Range("H132").Value = 9623868
$ws.Range("I132").Value = 12503754
$ws.Range("K132").Value = 37511262
$ws.Range("M132").Value = -37508732
